$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "47.340.40"
$ws.Range("E2").Value = "  +2.45%  "

# Row 3 - Ethereum
Set-TextValue "D3" "2.501.92"
$ws.Range("E3").Value = "  +2.07%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.11%  "

# Row 5 - BNB
Set-TextValue "D5" "324.01"
$ws.Range("E5").Value = "  +0.95%  "

# Row 6 - Solana
Set-TextValue "D6" "110.04"
$ws.Range("E6").Value = "  +4.92%  "

# Row 7 - XRP
$ws.Range("E7").Value = "  +1.14%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  +0.10%  "

# Row 10 - Avalanche
Set-TextValue "D10" "39.45"
$ws.Range("E10").Value = "  +9.71%  "

# Row 11 - Dogecoin
Set-TextValue "D11" "0.0815"
$ws.Range("E11").Value = "  +0.94%  "

# Row 12 - TRON
$ws.Range("E12").Value = "  +0.96%  "

# Row 13 - Chainlink
Set-TextValue "D13" "18.47"
$ws.Range("E13").Value = "  +1.01%  "

# Row 14 - Polkadot
$ws.Range("E14").Value = "  +1.83%  "

# Row 15 - WrappedliquidstakedEther2.0
Set-TextValue "D15" "2.894.52"
$ws.Range("E15").Value = "  +2.05%  "

# Row 16 - WrappedEther
Set-TextValue "D16" "2.498.81"
$ws.Range("E16").Value = "  +1.85%  "

# Row 17 - Polygon
Set-TextValue "D17" "0.859"
$ws.Range("E17").Value = "  +1.72%  "

# Row 18 - WrappedBTC
Set-TextValue "D18" "47.269.72"
$ws.Range("E18").Value = "  +2.65%  "

# Row 19 - InternetComputer(DFINITY)
Set-TextValue "D19" "12.86"
$ws.Range("E19").Value = "  +1.95%  "

# Row 20 - Uniswap
$ws.Range("E20").Value = "  +3.92%  "

# Row 21 - ShibaInu
Set-TextValue "D21" "0.0₃0941"
$ws.Range("E21").Value = "  +0.74%  "

# Row 22 - ImmutableX
Set-TextValue "D22" "2.72"
$ws.Range("E22").Value = "  +14.15%  "

# Row 23 - Litecoin
Set-TextValue "D23" "70.51"
$ws.Range("E23").Value = "  -0.98%  "

# Row 24 - BitcoinCash
Set-TextValue "D24" "248.18"
$ws.Range("E24").Value = "  +0.34%  "

# Row 25 - PancakeSwap
Set-TextValue "D25" "2.61"
$ws.Range("E25").Value = "  +3.81%  "

# Row 26 - EthereumClassic
Set-TextValue "D26" "26.09"
$ws.Range("E26").Value = "  +0.62%  "

# Row 27 - Dai
$ws.Range("E27").Value = "  -0.13%  "

# Row 28 & 29 - swap Toncoin / Cosmos
$ws.Range("B28").Value = "Cosmos"
$ws.Range("C28").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue "D28" "10.07"
$ws.Range("E28").Value = "  +3.90%  "

$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue "D29" "2.22"
$ws.Range("E29").Value = "  -3.20%  "

# Row 30 - InjectiveProtocol
Set-TextValue "D30" "35.37"
$ws.Range("E30").Value = "  +4.72%  "

# Row 31 - Kaspa
$ws.Range("E31").Value = "  +8.43%  "

# Row 32 - OKB
Set-TextValue "D32" "49.92"
$ws.Range("E32").Value = "  +1.22%  "

# Row 33 - Celestia
Set-TextValue "D33" "20.04"
$ws.Range("E33").Value = "  +0.46%  "

# Row 34 - Filecoin
$ws.Range("E34").Value = "  +1.99%  "

# Row 35 - Hedera
$ws.Range("E35").Value = "  +4.43%  "

# Row 36 - FirstDigitalUSD
$ws.Range("E36").Value = "  +0.33%  "

# Row 37 - ARBITRUM
$ws.Range("E37").Value = "  +5.06%  "

# Row 38 - RenderToken
$ws.Range("E38").Value = "  +3.94%  "

# Row 39 - LidoDAOToken
Set-TextValue "D39" "3.00"
$ws.Range("E39").Value = "  +1.99%  "

# Row 40 - Stellar
Set-TextValue "D40" "0.112"
$ws.Range("E40").Value = "  +1.38%  "

# Row 41 - Monero
Set-TextValue "D41" "121.79"
$ws.Range("E41").Value = "  -3.22%  "

# Row 42 - WEMIXToken
Set-TextValue "D42" "2.24"
$ws.Range("E42").Value = "  -0.55%  "

# Row 43 - EnergySwap
Set-TextValue "D43" "21.26"
$ws.Range("E43").Value = "  +1.51%  "

# Row 44 - VeChain
$ws.Range("E44").Value = "  +2.40%  "

# Row 45 - Maker
Set-TextValue "D45" "2.003.11"
$ws.Range("E45").Value = "  +1.73%  "

# Row 46 - NEARProtocol
$ws.Range("E46").Value = "  +4.15%  "

# Row 47 - ApeXProtocol
$ws.Range("E47").Value = "  -0.16%  "

# Row 48 - Stacks
$ws.Range("E48").Value = "  -4.24%  "

# Row 49 - FraxShare
$ws.Range("E49").Value = "  -0.66%  "

# Row 50 - THORChain
Set-TextValue "D50" "5.23"
$ws.Range("E50").Value = "  +4.04%  "

# Row 51 - MultiversX
Set-TextValue "D51" "56.78"
$ws.Range("E51").Value = "  +3.87%  "
